$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.953.58"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.554.04"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'206.73"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'21.61"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "'0.0859"
$ws.Range("D12").Value = "1.773.07"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.553.69"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "'3.71"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "26.952.02"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "'61.79"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "'214.50"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'7.26"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").Value = "'153.13"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "'6.66"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'14.88"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").Value = "1.376.56"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").Value = "'2.98"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "'0.972"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "'0.523"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'0.809"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "'0.993"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'5.48"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").Value = "'63.81"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").Value = "1.687.38"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'86.29"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "'0.0508"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  +0.35%  "
